$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 623.53845
$ws.Range("I19").Value = 554.2
$ws.Range("J19").Value = 666.875
$ws.Range("K19").Value = 554.2
$ws.Range("L19").Value = 666.875
$ws.Range("M19").Value = -379.2
$ws.Range("N19").Value = -1016.875
$ws.Range("H34").Value = 2785
$ws.Range("I34").Value = 1350.4166
$ws.Range("J34").Value = 20000
$ws.Range("K34").Value = 1350.4166
$ws.Range("L34").Value = 20000
$ws.Range("M34").Value = -1147.4166
$ws.Range("N34").Value = -20406
$ws.Range("H36").Value = 2785
$ws.Range("I36").Value = 1350.4166
$ws.Range("J36").Value = 20000
$ws.Range("K36").Value = 1350.4166
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = -635.4166
$ws.Range("N36").Value = -21430
$ws.Range("H43").Value = 2579.8
$ws.Range("I43").Value = 1950
$ws.Range("J43").Value = 2999.6667
$ws.Range("K43").Value = 1950
$ws.Range("L43").Value = 2999.6667
$ws.Range("M43").Value = -1881
$ws.Range("N43").Value = -3137.6667
$ws.Range("H53").Value = 45455590
$ws.Range("I53").Value = 71430090
$ws.Range("J53").Value = 235.25
$ws.Range("K53").Value = 71430090
$ws.Range("L53").Value = 235.25
$ws.Range("M53").Value = -71429453
$ws.Range("N53").Value = -1509.25
$ws.Range("H106").Value = 78434200
$ws.Range("I106").Value = 33336934
$ws.Range("K106").Value = 33336934
$ws.Range("M106").Value = -33336303
$ws.Range("H127").Value = 1604.6904
$ws.Range("I127").Value = 362.625
$ws.Range("J127").Value = 1896.9412
$ws.Range("K127").Value = 1087.875
$ws.Range("L127").Value = 5690.8236
$ws.Range("M127").Value = 3872.125
$ws.Range("N127").Value = -15610.8236
$ws.Range("H138").Value = 5316.4346
$ws.Range("I138").Value = 1042.6316
$ws.Range("J138").Value = 10555.29
$ws.Range("K138").Value = 3127.8948
$ws.Range("L138").Value = 31665.87
$ws.Range("M138").Value = 2012.1052
$ws.Range("N138").Value = -41945.87
$ws.Range("H141").Value = 1707
$ws.Range("I141").Value = 1707
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 5121
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 59
$ws.Range("N141").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1597.2174
$ws.Range("I2").Value = 1150.8572
$ws.Range("J2").Value = 2291.5557
$ws.Range("K2").Value = 1150.8572
$ws.Range("L2").Value = 2291.5557
$ws.Range("M2").Value = -1037.8572
$ws.Range("N2").Value = -2517.5557
$ws.Range("H24").Value = 31564
$ws.Range("J24").Value = 31564
$ws.Range("L24").Value = 31564
$ws.Range("N24").Value = -32312
$ws.Range("H45").Value = 3683.853
$ws.Range("I45").Value = 4146.448
$ws.Range("J45").Value = 1000.8
$ws.Range("K45").Value = 4146.448
$ws.Range("L45").Value = 1000.8
$ws.Range("M45").Value = -3769.448
$ws.Range("N45").Value = -1754.8
$ws.Range("H61").Value = 2789.2754
$ws.Range("I61").Value = 2971.6724
$ws.Range("J61").Value = 1827.5454
$ws.Range("K61").Value = 2971.6724
$ws.Range("L61").Value = 1827.5454
$ws.Range("M61").Value = -2759.6724
$ws.Range("N61").Value = -2251.5454
$ws.Range("H100").Value = 31564
$ws.Range("J100").Value = 31564
$ws.Range("L100").Value = 31564
$ws.Range("N100").Value = -33728
$ws.Range("H116").Value = 1597.2174
$ws.Range("I116").Value = 1150.8572
$ws.Range("J116").Value = 2291.5557
$ws.Range("K116").Value = 1150.8572
$ws.Range("L116").Value = 2291.5557
$ws.Range("M116").Value = 1143.1428
$ws.Range("N116").Value = -6879.5557
$ws.Range("H122").Value = 6412759
$ws.Range("I122").Value = 6412759
$ws.Range("K122").Value = 19238277
$ws.Range("M122").Value = -19235827
$ws.Range("H136").Value = 2789.2754
$ws.Range("I136").Value = 2971.6724
$ws.Range("J136").Value = 1827.5454
$ws.Range("K136").Value = 8915.0172
$ws.Range("L136").Value = 5482.6362
$ws.Range("M136").Value = -6365.0172
$ws.Range("N136").Value = -10582.6362
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1597.2174
$ws.Range("I3").Value = 1150.8572
$ws.Range("J3").Value = 2291.5557
$ws.Range("K3").Value = 1150.8572
$ws.Range("L3").Value = 2291.5557
$ws.Range("M3").Value = -1036.8572
$ws.Range("N3").Value = -2519.5557
$ws.Range("H105").Value = 2387.1428
$ws.Range("I105").Value = 1822
$ws.Range("J105").Value = 3800
$ws.Range("K105").Value = 1822
$ws.Range("L105").Value = 3800
$ws.Range("M105").Value = -75
$ws.Range("N105").Value = -7294
$ws.Range("H107").Value = 805.1875
$ws.Range("I107").Value = 784.0769
$ws.Range("J107").Value = 896.6667
$ws.Range("K107").Value = 784.0769
$ws.Range("L107").Value = 896.6667
$ws.Range("M107").Value = 1135.9231
$ws.Range("N107").Value = -4736.6667
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5171.7637
$ws.Range("I31").Value = 1741.8
$ws.Range("K31").Value = 1741.8
$ws.Range("M31").Value = -1446.8
$ws.Range("H34").Value = 5171.7637
$ws.Range("I34").Value = 1741.8
$ws.Range("K34").Value = 1741.8
$ws.Range("M34").Value = -1539.8
$ws.Range("H99").Value = 5960238.5
$ws.Range("J99").Value = 11369600
$ws.Range("L99").Value = 11369600
$ws.Range("N99").Value = -11372596
$ws.Range("H122").Value = 2684.7896
$ws.Range("I122").Value = 2094.75
$ws.Range("J122").Value = 3696.2856
$ws.Range("K122").Value = 6284.25
$ws.Range("L122").Value = 11088.8568
$ws.Range("M122").Value = -3834.25
$ws.Range("N122").Value = -15988.8568
$ws.Range("H126").Value = 5960238.5
$ws.Range("J126").Value = 11369600
$ws.Range("L126").Value = 34108800
$ws.Range("N126").Value = -34113740
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 3406
$ws.Range("I63").Value = 312
$ws.Range("K63").Value = 936
$ws.Range("M63").Value = -187
$ws.Range("H64").Value = 4137.778
$ws.Range("I64").Value = 4953
$ws.Range("K64").Value = 14859
$ws.Range("M64").Value = -14589
$ws.Range("H66").Value = 3406
$ws.Range("I66").Value = 312
$ws.Range("K66").Value = 2808
$ws.Range("M66").Value = 936
$ws.Range("H67").Value = 4137.778
$ws.Range("I67").Value = 4953
$ws.Range("K67").Value = 14859
$ws.Range("M67").Value = -13923
$ws.Range("H70").Value = 1646.7693
$ws.Range("I70").Value = 656.44446
$ws.Range("J70").Value = 3875
$ws.Range("K70").Value = 1969.33338
$ws.Range("L70").Value = 11625
$ws.Range("M70").Value = -1654.33338
$ws.Range("N70").Value = -12255
$ws.Range("H73").Value = 1646.7693
$ws.Range("I73").Value = 656.44446
$ws.Range("J73").Value = 3875
$ws.Range("K73").Value = 1969.33338
$ws.Range("L73").Value = 11625
$ws.Range("M73").Value = -877.33338
$ws.Range("N73").Value = -13809
$ws.Range("H132").Value = 1612.1111
$ws.Range("I132").Value = 1084
$ws.Range("J132").Value = 2668.3333
$ws.Range("K132").Value = 9756
$ws.Range("L132").Value = 24014.9997
$ws.Range("M132").Value = -7226
$ws.Range("N132").Value = -29074.9997
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 47620330
$ws.Range("I113").Value = 83334376
$ws.Range("J113").Value = 1599.2222
$ws.Range("K113").Value = 83334376
$ws.Range("L113").Value = 1599.2222
$ws.Range("M113").Value = -83332206
$ws.Range("N113").Value = -5939.2222
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3499.9167
$ws.Range("I61").Value = 2542.7144
$ws.Range("J61").Value = 4840
$ws.Range("K61").Value = 2542.7144
$ws.Range("L61").Value = 4840
$ws.Range("M61").Value = -2340.7144
$ws.Range("N61").Value = -5244
$ws.Range("H113").Value = 3499.9167
$ws.Range("I113").Value = 2542.7144
$ws.Range("J113").Value = 4840
$ws.Range("K113").Value = 2542.7144
$ws.Range("L113").Value = 4840
$ws.Range("M113").Value = -372.7143999999998
$ws.Range("N113").Value = -9180
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 83333510
$ws.Range("I107").Value = 125000140
$ws.Range("K107").Value = 375000420
$ws.Range("M107").Value = -374998500
$ws.Range("H122").Value = 1661.3
$ws.Range("I122").Value = 1825.2
$ws.Range("J122").Value = 1169.6
$ws.Range("K122").Value = 5475.6
$ws.Range("L122").Value = 3508.8
$ws.Range("M122").Value = -3025.6
$ws.Range("N122").Value = -8408.8
